$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 3993.3333
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 4990
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 4990
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -5214

$ws.Range("H14").Value = 3993.3333
$ws.Range("I14").Value = 2000
$ws.Range("J14").Value = 4990
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 4990
$ws.Range("M14").Value = -1809
$ws.Range("N14").Value = -5372

$ws.Range("H15").Value = 2521.514
$ws.Range("I15").Value = 2521.514
$ws.Range("K15").Value = 7564.542
$ws.Range("M15").Value = -7395.542

$ws.Range("H33").Value = 770226
$ws.Range("I33").Value = 1001064.4
$ws.Range("J33").Value = 764.6667
$ws.Range("K33").Value = 1001064.4
$ws.Range("L33").Value = 764.6667
$ws.Range("M33").Value = -1000835.4
$ws.Range("N33").Value = -1222.6667

$ws.Range("H38").Value = 1295.56
$ws.Range("I38").Value = 159.26666
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 477.79998
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = -105.79998
$ws.Range("N38").Value = -9744

$ws.Range("H39").Value = 393.2
$ws.Range("I39").Value = 80.55
$ws.Range("J39").Value = 705.85
$ws.Range("K39").Value = 241.65
$ws.Range("L39").Value = 2117.55
$ws.Range("M39").Value = 54.35000000000002
$ws.Range("N39").Value = -2709.55

$ws.Range("H98").Value = 2450.7144
$ws.Range("I98").Value = 1925.8334
$ws.Range("J98").Value = 5600
$ws.Range("K98").Value = 1925.8334
$ws.Range("L98").Value = 5600
$ws.Range("M98").Value = -427.8334
$ws.Range("N98").Value = -8596

$ws.Range("H116").Value = 2817.8857
$ws.Range("I116").Value = 1909.68
$ws.Range("J116").Value = 5088.4
$ws.Range("K116").Value = 1909.68
$ws.Range("L116").Value = 5088.4
$ws.Range("M116").Value = 1532.32
$ws.Range("N116").Value = -11972.4

$ws.Range("H122").Value = 2450.7144
$ws.Range("I122").Value = 1925.8334
$ws.Range("J122").Value = 5600
$ws.Range("K122").Value = 5777.5002
$ws.Range("L122").Value = 16800
$ws.Range("M122").Value = -3327.5002
$ws.Range("N122").Value = -21700

$ws.Range("H132").Value = 8514.194
$ws.Range("I132").Value = 8957.950000000001
$ws.Range("J132").Value = 7959.5
$ws.Range("K132").Value = 26873.85
$ws.Range("L132").Value = 23878.5
$ws.Range("M132").Value = -24343.85
$ws.Range("N132").Value = -28938.5

$ws.Range("H137").Value = 9435765
$ws.Range("I137").Value = 15153427
$ws.Range("J137").Value = 1622.75
$ws.Range("K137").Value = 45460281
$ws.Range("L137").Value = 4868.25
$ws.Range("M137").Value = -45457731
$ws.Range("N137").Value = -9968.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3978457.8
$ws.Range("I32").Value = 5805.5
$ws.Range("J32").Value = 33376086
$ws.Range("K32").Value = 5805.5
$ws.Range("L32").Value = 33376086
$ws.Range("M32").Value = -5518.5
$ws.Range("N32").Value = -33376660

$ws.Range("H61").Value = 2307.6052
$ws.Range("I61").Value = 1442.3478
$ws.Range("J61").Value = 3634.3333
$ws.Range("K61").Value = 1442.3478
$ws.Range("L61").Value = 3634.3333
$ws.Range("M61").Value = -1230.3478
$ws.Range("N61").Value = -4058.3333

$ws.Range("H122").Value = 1528.1132
$ws.Range("I122").Value = 1089.6511
$ws.Range("K122").Value = 3268.9533
$ws.Range("M122").Value = -818.9533000000001

$ws.Range("H132").Value = 2362357.5
$ws.Range("I132").Value = 1699.8182
$ws.Range("J132").Value = 6257442.5
$ws.Range("K132").Value = 5099.4546
$ws.Range("L132").Value = 18772327.5
$ws.Range("M132").Value = -2569.4546
$ws.Range("N132").Value = -18777387.5

$ws.Range("H136").Value = 2307.6052
$ws.Range("I136").Value = 1442.3478
$ws.Range("J136").Value = 3634.3333
$ws.Range("K136").Value = 4327.0434
$ws.Range("L136").Value = 10902.9999
$ws.Range("M136").Value = -1777.0434
$ws.Range("N136").Value = -16002.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1749.5667
$ws.Range("I20").Value = 1384.2858
$ws.Range("J20").Value = 2601.889
$ws.Range("K20").Value = 1384.2858
$ws.Range("L20").Value = 2601.889
$ws.Range("M20").Value = -1137.2858
$ws.Range("N20").Value = -3095.889

$ws.Range("H99").Value = 2148.4814
$ws.Range("I99").Value = 1911.6875
$ws.Range("J99").Value = 2492.9092
$ws.Range("K99").Value = 1911.6875
$ws.Range("L99").Value = 2492.9092
$ws.Range("M99").Value = -413.6875
$ws.Range("N99").Value = -5488.9092

$ws.Range("H134").Value = 4833.021
$ws.Range("I134").Value = 2144.32
$ws.Range("J134").Value = 7755.522
$ws.Range("K134").Value = 6432.960000000001
$ws.Range("L134").Value = 23266.566
$ws.Range("M134").Value = -3897.960000000001
$ws.Range("N134").Value = -28336.566

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6287.8237
$ws.Range("I7").Value = 10051.5
$ws.Range("J7").Value = 911.1429000000001
$ws.Range("K7").Value = 10051.5
$ws.Range("L7").Value = 911.1429000000001
$ws.Range("M7").Value = -9938.5
$ws.Range("N7").Value = -1137.1429

$ws.Range("H10").Value = 1990
$ws.Range("I10").Value = 1990
$ws.Range("K10").Value = 1990
$ws.Range("M10").Value = -1851

$ws.Range("H94").Value = 1333.6471
$ws.Range("I94").Value = 1262
$ws.Range("J94").Value = 1343.2
$ws.Range("K94").Value = 1262
$ws.Range("L94").Value = 1343.2
$ws.Range("M94").Value = -811
$ws.Range("N94").Value = -2245.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 647.5349
$ws.Range("I5").Value = 420.625
$ws.Range("J5").Value = 934.1579
$ws.Range("K5").Value = 1261.875
$ws.Range("L5").Value = 2802.4737
$ws.Range("M5").Value = -1149.875
$ws.Range("N5").Value = -3026.4737

$ws.Range("H8").Value = 73.5625
$ws.Range("I8").Value = 73.5625
$ws.Range("K8").Value = 220.6875
$ws.Range("M8").Value = -81.6875

$ws.Range("H68").Value = 812976.3
$ws.Range("I68").Value = 1042.5385
$ws.Range("J68").Value = 1058444.6
$ws.Range("K68").Value = 3127.6155
$ws.Range("L68").Value = 3175333.8
$ws.Range("M68").Value = -2316.6155
$ws.Range("N68").Value = -3176955.8

$ws.Range("H71").Value = 812976.3
$ws.Range("I71").Value = 1042.5385
$ws.Range("J71").Value = 1058444.6
$ws.Range("K71").Value = 9382.846500000001
$ws.Range("L71").Value = 9526001.4
$ws.Range("M71").Value = -5326.846500000001
$ws.Range("N71").Value = -9534113.4

$ws.Range("H113").Value = 490.42856
$ws.Range("I113").Value = 468.77777
$ws.Range("J113").Value = 513.35297
$ws.Range("K113").Value = 1406.33331
$ws.Range("L113").Value = 1540.05891
$ws.Range("M113").Value = 763.66669
$ws.Range("N113").Value = -5880.05891

$ws.Range("H135").Value = 647.5349
$ws.Range("I135").Value = 420.625
$ws.Range("J135").Value = 934.1579
$ws.Range("K135").Value = 3785.625
$ws.Range("L135").Value = 8407.4211
$ws.Range("M135").Value = -1250.625
$ws.Range("N135").Value = -13477.4211

$ws.Range("H141").Value = 3203.75
$ws.Range("I141").Value = 2026
$ws.Range("J141").Value = 5166.6665
$ws.Range("K141").Value = 6078
$ws.Range("L141").Value = 15499.9995
$ws.Range("M141").Value = -898
$ws.Range("N141").Value = -25859.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7828.5713
$ws.Range("I70").Value = 6075
$ws.Range("K70").Value = 6075
$ws.Range("M70").Value = -5805

$ws.Range("H73").Value = 7828.5713
$ws.Range("I73").Value = 6075
$ws.Range("K73").Value = 6075
$ws.Range("M73").Value = -5139

$ws.Range("H113").Value = 1486.3334
$ws.Range("I113").Value = 1319.2727
$ws.Range("J113").Value = 1945.75
$ws.Range("K113").Value = 1319.2727
$ws.Range("L113").Value = 1945.75
$ws.Range("M113").Value = 850.7273
$ws.Range("N113").Value = -6285.75

$ws.Range("H132").Value = 2355.8965
$ws.Range("I132").Value = 2923.3
$ws.Range("J132").Value = 2057.2632
$ws.Range("K132").Value = 8769.900000000001
$ws.Range("L132").Value = 6171.7896
$ws.Range("M132").Value = -6239.900000000001
$ws.Range("N132").Value = -11231.7896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 251860.7
$ws.Range("I93").Value = 386229.22
$ws.Range("J93").Value = 2319.1428
$ws.Range("K93").Value = 386229.22
$ws.Range("L93").Value = 2319.1428
$ws.Range("M93").Value = -384981.22
$ws.Range("N93").Value = -4815.1428

$ws.Range("H132").Value = 29445574
$ws.Range("I132").Value = 38504780
$ws.Range("J132").Value = 3160.375
$ws.Range("K132").Value = 115514340
$ws.Range("L132").Value = 9481.125
$ws.Range("M132").Value = -115511810
$ws.Range("N132").Value = -14541.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 915.65515
$ws.Range("I113").Value = 470.7143
$ws.Range("J113").Value = 1057.2273
$ws.Range("K113").Value = 1412.1429
$ws.Range("L113").Value = 3171.6819
$ws.Range("M113").Value = 757.8571000000002
$ws.Range("N113").Value = -7511.6819

$ws.Range("H132").Value = 2435.946
$ws.Range("I132").Value = 2134.111
$ws.Range("J132").Value = 2721.8948
$ws.Range("K132").Value = 6402.333
$ws.Range("L132").Value = 8165.6844
$ws.Range("M132").Value = -3872.333
$ws.Range("N132").Value = -13225.6844
